$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5: "sha_url" / "http://www.sha1-online.com/" -----------------

# A5 ("sha_url"): same look as the other "Key" column cells (A1:A4) -
# copy A4's format first so it reuses the existing cell style, then set
# the text.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = "sha_url"

# B5 ("http://www.sha1-online.com/"): value cell, styled like a hyperlink -
# underlined Calibri, black, left aligned (new font/style vs. the
# Aptos-Narrow one used by B2:B4).
$ws.Range("B5").Value = "http://www.sha1-online.com/"
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("B5").Font.Name = "Calibri"
$ws.Range("B5").Font.Underline = $true
$ws.Range("B5").Font.Color = 0

# Row 5 height, matching the new row's taller height.
$ws.Rows.Item(5).RowHeight = 19.5

# Column A widens slightly to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 18.57
